$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Redactor text" -> three runs "Redactor" / " " / "text"
#    (the middle space loses its w:lang="en-US" mark, matching the captured
#    edit where the space was retyped without the language tag).
# ---------------------------------------------------------------------------
$titleRange = $d.Content.Duplicate
$titleRange.Find.Execute("Redactor text", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr><w:t>Redactor</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr><w:t>text</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$titleRange.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2) Remove the "1254" run entirely, leaving its paragraph empty.
# ---------------------------------------------------------------------------
$numRange = $d.Content.Duplicate
$numRange.Find.Execute("1254", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$numRange.Text = ""

# ---------------------------------------------------------------------------
# 3) Mark the run that holds the inline picture as "do not spell-check"
#    (adds <w:noProof/> to that run's rPr).
# ---------------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $picRange = $d.InlineShapes.Item(1).Range
    $picRange.NoProofing = $true
}

# ---------------------------------------------------------------------------
# 4) City name swap: "Казань" -> "Россия" (same run formatting, text only).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Казань", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Россия", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Style "Default Paragraph Font" becomes semi-hidden in the style gallery
#    (w:semiHidden on styleId "a0"). Guarded in case the host's Style object
#    does not expose a writable Hidden property.
# ---------------------------------------------------------------------------
$defaultParaFont = $d.Styles("Default Paragraph Font")
try {
    $defaultParaFont.Hidden = $true
} catch {
    Write-Output ("Style.Hidden not settable on this host: " + $_.Exception.Message)
}
